$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The fertilizer abbreviation "KAS" (Kalkammonsalpeter) used in the
# "Notice" column of the nitrogen application tables is renamed to its
# English equivalent "CAN" (Calcium Ammonium Nitrate) for consistency.
$cells = @("E44","E45","E55","E56","E57","E66","E67","E68","E77","E78","E79")
foreach ($c in $cells) {
    $ws.Range($c).Value = "CAN"
}

# Mirror the last-selected cell recorded in the saved workbook.
$ws.Range("E79").Select()
